$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: A2 changes from 0 to 1 (B2 stays 32)
$ws.Range("A2").Value = 1

# Update row 3: A3 changes from 1 to 0, B3 changes from 14 to 25
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 25

# Remove row 4 entirely (A4=2, B4=11 are deleted)
$ws.Rows.Item(4).Delete()
